$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)   # quality_comparison
$ws2 = $wb.Worksheets.Item(2)   # computational_comparison

# --- Build the two new border-only cell styles on sheet1, then propagate
#     them to the equivalent cells on sheet2 via copy/paste-formats so that
#     no extra (unused) style entries get created in the shared style table.

# Style A: top + bottom thin border only (no font/alignment changes)
$c1 = $ws1.Range("C1")
$c1.Style = "Normal"
$c1.Borders.Item(8).LineStyle = 1   # xlEdgeTop
$c1.Borders.Item(9).LineStyle = 1   # xlEdgeBottom

# Style B: right + top + bottom thin border only
$d1 = $ws1.Range("D1")
$d1.Style = "Normal"
$d1.Borders.Item(10).LineStyle = 1  # xlEdgeRight
$d1.Borders.Item(8).LineStyle = 1   # xlEdgeTop
$d1.Borders.Item(9).LineStyle = 1   # xlEdgeBottom

# Propagate Style A to the matching cells of the computational_comparison sheet
$c1.Copy()
$ws2.Range("C1").PasteSpecial(-4122)   # xlPasteFormats
$ws2.Range("F1").PasteSpecial(-4122)   # xlPasteFormats

# Propagate Style B to the matching cells of the computational_comparison sheet
$d1.Copy()
$ws2.Range("D1").PasteSpecial(-4122)   # xlPasteFormats
$ws2.Range("G1").PasteSpecial(-4122)   # xlPasteFormats

# --- Anonymize "fedcore" -> "approach"
$ws1.Range("C2").Value = "approach"
$ws2.Range("C2").Value = "approach"
$ws2.Range("F2").Value = "approach"

# --- Drop the stray empty inline-string cell G5 on computational_comparison
$ws2.Range("G5").ClearContents()
